$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 201.75
$ws.Range("I2").Value = 133
$ws.Range("K2").Value = 133
$ws.Range("M2").Value = -20
$ws.Range("H111").Value = 1673.8334
$ws.Range("I111").Value = 1398.375
$ws.Range("J111").Value = 2224.75
$ws.Range("K111").Value = 4195.125
$ws.Range("L111").Value = 6674.25
$ws.Range("M111").Value = -1128.125
$ws.Range("N111").Value = -12808.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1594.1428
$ws.Range("I2").Value = 191
$ws.Range("K2").Value = 191
$ws.Range("M2").Value = -78
$ws.Range("H116").Value = 1594.1428
$ws.Range("I116").Value = 191
$ws.Range("K116").Value = 191
$ws.Range("M116").Value = 2103

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1594.1428
$ws.Range("I3").Value = 191
$ws.Range("K3").Value = 191
$ws.Range("M3").Value = -77
$ws.Range("H94").Value = 2718.818
$ws.Range("I94").Value = 981.8
$ws.Range("J94").Value = 4166.3335
$ws.Range("K94").Value = 981.8
$ws.Range("L94").Value = 4166.3335
$ws.Range("M94").Value = -530.8
$ws.Range("N94").Value = -5068.3335

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1424.25
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1424.25
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1424.25
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2124.25
$ws.Range("H58").Value = 6486.5386
$ws.Range("I58").Value = 7228.8
$ws.Range("K58").Value = 7228.8
$ws.Range("M58").Value = -7025.8
$ws.Range("H99").Value = 6636.8887
$ws.Range("I99").Value = 5980.077
$ws.Range("K99").Value = 5980.077
$ws.Range("M99").Value = -4482.077
$ws.Range("H126").Value = 6636.8887
$ws.Range("I126").Value = 5980.077
$ws.Range("K126").Value = 17940.231
$ws.Range("M126").Value = -15470.231
$ws.Range("H132").Value = 2426.75
$ws.Range("I132").Value = 2426.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7280.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4750.25
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1702.45
$ws.Range("I134").Value = 1225.4445
$ws.Range("K134").Value = 3676.3335
$ws.Range("M134").Value = -1141.3335
$ws.Range("H136").Value = 6486.5386
$ws.Range("I136").Value = 7228.8
$ws.Range("K136").Value = 21686.4
$ws.Range("M136").Value = -19136.4
$ws.Range("H141").Value = 256952
$ws.Range("J141").Value = 279946.66
$ws.Range("L141").Value = 279946.66
$ws.Range("N141").Value = -290306.66

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 804
$ws.Range("I18").Value = 766.25
$ws.Range("J18").Value = 955
$ws.Range("K18").Value = 2298.75
$ws.Range("L18").Value = 2865
$ws.Range("M18").Value = -2129.75
$ws.Range("N18").Value = -3203
$ws.Range("H46").Value = 95.666664
$ws.Range("I46").Value = 95.666664
$ws.Range("K46").Value = 286.999992
$ws.Range("M46").Value = -195.999992
$ws.Range("H51").Value = 3066
$ws.Range("I51").Value = 1198
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 3594
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = -3134
$ws.Range("N51").Value = -12920
$ws.Range("H80").Value = 5786
$ws.Range("J80").Value = 10500
$ws.Range("L80").Value = 31500
$ws.Range("N80").Value = -33372
$ws.Range("H83").Value = 5786
$ws.Range("J83").Value = 10500
$ws.Range("L83").Value = 94500
$ws.Range("N83").Value = -103860
$ws.Range("H98").Value = 584.36365
$ws.Range("I98").Value = 535.5
$ws.Range("J98").Value = 612.2857
$ws.Range("K98").Value = 1606.5
$ws.Range("L98").Value = 1836.8571
$ws.Range("M98").Value = -108.5
$ws.Range("N98").Value = -4832.8571
$ws.Range("H107").Value = 1603.875
$ws.Range("J107").Value = 1865.421
$ws.Range("L107").Value = 5596.263
$ws.Range("N107").Value = -9436.262999999999
$ws.Range("H113").Value = 745.25
$ws.Range("J113").Value = 745.25
$ws.Range("L113").Value = 2235.75
$ws.Range("N113").Value = -6575.75
$ws.Range("H122").Value = 1029.0714
$ws.Range("I122").Value = 438.2
$ws.Range("J122").Value = 1357.3334
$ws.Range("K122").Value = 3943.8
$ws.Range("L122").Value = 12216.0006
$ws.Range("M122").Value = -1493.8
$ws.Range("N122").Value = -17116.0006
$ws.Range("H140").Value = 1600

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H6").Value = 1752.5
$ws.Range("J6").Value = 2005
$ws.Range("L6").Value = 2005
$ws.Range("N6").Value = -2231
$ws.Range("H16").Value = 1752.5
$ws.Range("J16").Value = 2005
$ws.Range("L16").Value = 2005
$ws.Range("N16").Value = -2505
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 42995
$ws.Range("I62").Value = 42995
$ws.Range("K62").Value = 42995
$ws.Range("M62").Value = -42309
$ws.Range("H65").Value = 42995
$ws.Range("I65").Value = 42995
$ws.Range("K65").Value = 128985
$ws.Range("M65").Value = -125553
$ws.Range("H102").Value = 938
$ws.Range("I102").Value = 763.6667
$ws.Range("K102").Value = 763.6667
$ws.Range("M102").Value = 858.3333
$ws.Range("H105").Value = 9000
$ws.Range("J105").Value = 9000
$ws.Range("L105").Value = 9000
$ws.Range("N105").Value = -15988

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9449.666999999999
$ws.Range("J40").Value = 23999.5
$ws.Range("L40").Value = 23999.5
$ws.Range("N40").Value = -24271.5
$ws.Range("H46").Value = 2337.1177
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2337.1177
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2337.1177
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2713.1177
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H63").Value = 40445
$ws.Range("I63").Value = 40445
$ws.Range("K63").Value = 40445
$ws.Range("M63").Value = -39696
$ws.Range("H66").Value = 40445
$ws.Range("I66").Value = 40445
$ws.Range("K66").Value = 121335
$ws.Range("M66").Value = -117591
$ws.Range("H100").Value = 1470.5714
$ws.Range("I100").Value = 1432.6666
$ws.Range("J100").Value = 1499
$ws.Range("K100").Value = 1432.6666
$ws.Range("L100").Value = 1499
$ws.Range("M100").Value = -891.6666
$ws.Range("N100").Value = -2581
$ws.Range("H136").Value = 4738.5
$ws.Range("I136").Value = 4435.625
$ws.Range("K136").Value = 13306.875
$ws.Range("M136").Value = -10756.875
